# Weekly price update: a new week's price record for "Ajo" (Chino / Primera)
# at "Terminal La Palmera de La Serena" is inserted at the top of its block
# (row 143), pushing all subsequent rows of that sheet down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 143 - shifts the existing rows 143:165 down
# to 144:166 and grows the used range to A1:R166.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row with the new week's data.
$ws.Cells.Item(143, 1).Value = 8
$ws.Cells.Item(143, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(143, 3).Value = "Coquimbo"
$ws.Cells.Item(143, 4).Value = 44505
$ws.Cells.Item(143, 5).Value = 4
$ws.Cells.Item(143, 6).Value = 100112003
$ws.Cells.Item(143, 7).Value = "Ajo"
$ws.Cells.Item(143, 8).Value = "Chino"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 660
$ws.Cells.Item(143, 11).Value = 19000
$ws.Cells.Item(143, 12).Value = 20000
$ws.Cells.Item(143, 13).Value = 19500
$ws.Cells.Item(143, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(143, 15).Value = "China"
$ws.Cells.Item(143, 16).Value = 1950
$ws.Cells.Item(143, 17).Value = 10
$ws.Cells.Item(143, 18).Value = "Hortaliza"
